# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"     = @{ "F3" = 8516; "F4" = 6220; "F5" = 544; "F6" = 116; "F9" = 331; "F10" = 1202 }
    "全部类型" = @{ "F3" = 8516; "F4" = 6220; "F5" = 544; "F6" = 116; "F9" = 331; "F14" = 1202 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
